$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original data set (rows 64-118) is a time series of weekly price
# observations. This edit adds four new weekly observations, inserted at
# their chronological slot among the existing rows (pushing the rows that
# follow each insertion point down by one). Insert ascending so each
# insertion position already accounts for the rows added before it.
$ws.Rows(64).Insert()
$ws.Rows(102).Insert()
$ws.Rows(109).Insert()
$ws.Rows(113).Insert()

function Set-PriceRow {
    param($Row, $Fecha, $Volumen, $PrecioMinimo, $PrecioMaximo, $PrecioPromedio, $Unidad, $Origen, $PrecioKg, $KgUnidades)

    $ws.Cells.Item($Row, 1).Value = 7
    $ws.Cells.Item($Row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
    $ws.Cells.Item($Row, 3).Value = "Ñuble"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = 16
    $ws.Cells.Item($Row, 6).Value = 100112032
    $ws.Cells.Item($Row, 7).Value = "Zapallo italiano"
    $ws.Cells.Item($Row, 8).Value = "Sin especificar"
    $ws.Cells.Item($Row, 9).Value = "Primera"
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMinimo
    $ws.Cells.Item($Row, 12).Value = $PrecioMaximo
    $ws.Cells.Item($Row, 13).Value = $PrecioPromedio
    $ws.Cells.Item($Row, 14).Value = $Unidad
    $ws.Cells.Item($Row, 15).Value = $Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = $KgUnidades
    $ws.Cells.Item($Row, 18).Value = "Hortaliza"
}

Set-PriceRow 64  44447 160 16000 17000 16500 "$/caja 50 unidades" "Región de Arica y Parinacota" 330 50
Set-PriceRow 102 44441 160 16000 17000 16500 "$/caja 50 unidades" "Región de Arica y Parinacota" 330 50
Set-PriceRow 109 44446 160 16000 17000 16500 "$/caja 50 unidades" "Región de Arica y Parinacota" 330 50
Set-PriceRow 113 44442 120 16000 17000 16500 "$/caja 50 unidades" "Región de Arica y Parinacota" 330 50

# Make sure the date cells keep the same date number format as the rest of
# column D (Insert() already copies the format from the row above, but set
# it explicitly too so the new rows are not missing it for any reason).
$ws.Range("D64").NumberFormat = $ws.Range("D65").NumberFormat
$ws.Range("D102").NumberFormat = $ws.Range("D101").NumberFormat
$ws.Range("D109").NumberFormat = $ws.Range("D108").NumberFormat
$ws.Range("D113").NumberFormat = $ws.Range("D112").NumberFormat
